$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the three new columns, re-using the same style as the other
# header cells (F1:H1 get the header style copied from E1).
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Boolean outlier flags for rows 2-18. Row 5 is flagged TRUE for all three
# algorithms; every other row is FALSE.
$values = @{
    2  = $false
    3  = $false
    4  = $false
    5  = $true
    6  = $false
    7  = $false
    8  = $false
    9  = $false
    10 = $false
    11 = $false
    12 = $false
    13 = $false
    14 = $false
    15 = $false
    16 = $false
    17 = $false
    18 = $false
}

foreach ($row in $values.Keys) {
    $val = $values[$row]
    $ws.Range("F$row").Value = $val
    $ws.Range("G$row").Value = $val
    $ws.Range("H$row").Value = $val
}
